{"js": "// Update bibliography entries: expand abbreviated journal names to their\n// full titles and drop the trailing \"doi:\"/URL hyperlinks (per the\n// \"added updates for phrasing prior to submission\" commit).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Helper: find the first match of `needle` inside a paragraph and replace\n// it with `replacement`, preserving the run's existing formatting.\nasync function replaceInParagraph(paragraph, needle, replacement) {\n  const results = paragraph.getRange().search(needle, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + needle);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Helper: find the first match of `needle` inside a paragraph and clear it\n// (used to remove a doi/URL run together with its enclosing hyperlink).\nasync function removeInParagraph(paragraph, needle) {\n  const results = paragraph.getRange().search(needle, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + needle);\n  }\n  results.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Reference 1 (paragraph 90): Mangtani et al., Clin Infect Dis 2014 ---\nconst p90 = paragraphs.items[90];\nawait replaceInParagraph(\n  p90,\n  \"Clin Infect Dis\",\n  \"Clinical infectious diseases : an official publication of the Infectious Diseases Society of America\"\n);\nawait replaceInParagraph(p90, \":470\\u201380. doi:\", \":470\\u201380.\");\nawait removeInParagraph(p90, \"10.1093/cid/cit790\");\n\n// --- Reference 2 (paragraph 91): Abubakar et al., Health Technol Assess 2013 ---\nconst p91 = paragraphs.items[91];\nawait replaceInParagraph(p91, \"Health Technol Assess\", \"Health technology assessment\");\nawait replaceInParagraph(p91, \":1\\u2013372, v\\u2013vi. doi:\", \":1\\u2013372, v\\u2013vi.\");\nawait removeInParagraph(p91, \"10.3310/hta17370\");\n\n// --- Reference 3 (paragraph 92): Zwerling et al., PLoS Med 2011 ---\nconst p92 = paragraphs.items[92];\nawait replaceInParagraph(p92, \"PLoS Med\", \"PLoS medicine\");\nawait replaceInParagraph(p92, \". doi:\", \":e1001012.\");\nawait removeInParagraph(p92, \"10.1371/journal.pmed.1001012\");\n\n// --- Reference 4 (paragraph 93): Rodrigues et al., Int J Epidemiol 1993 ---\nconst p93 = paragraphs.items[93];\nawait replaceInParagraph(p93, \"Int J Epidemiol\", \"International journal of epidemiology\");\nawait removeInParagraph(p93, \"http://www.ncbi.nlm.nih.gov/pubmed/8144299\");\n\n// --- Reference 5 (paragraph 94): Colditz et al., JAMA 1994 ---\nconst p94 = paragraphs.items[94];\nawait replaceInParagraph(p94, \":698. doi:\", \":698.\");\nawait removeInParagraph(p94, \"10.1001/jama.1994.03510330076038\");\n\n// --- Reference 6 (paragraph 95): Mangtani et al., Clin Infect Dis 2014 (dup) ---\nconst p95 = paragraphs.items[95];\nawait replaceInParagraph(p95, \"Clin Infect Dis\", \"Clinical Infectious Diseases\");\nawait replaceInParagraph(p95, \":470\\u201380. doi:\", \":470\\u201380.\");\nawait removeInParagraph(p95, \"10.1093/cid/cit790\");\n\n// --- Reference 7 (paragraph 96): Hart & Sutherland, Bull World Health Organ 1972 ---\nconst p96 = paragraphs.items[96];\nawait replaceInParagraph(p96, \"Bull World Health Organ\", \"The American Statistician\");\nawait replaceInParagraph(p96, \":371\\u201385. doi:\", \":371\\u201385.\");\nawait removeInParagraph(p96, \"10.1136/bmj.2.6082.293\");\n\n// --- Reference 10 (paragraph 99): Stevenson et al., epiR package 2017 ---\nconst p99 = paragraphs.items[99];\nconst epiRTitle = p99\n  .getRange()\n  .search(\"epiR: Tools for the Analysis of Epidemiological Data\", { matchCase: true });\nepiRTitle.load(\"text,font\");\nawait context.sync();\nepiRTitle.items[0].font.set({ italic: true });\nawait context.sync();\n\nawait removeInParagraph(p99, \"https://cran.r-project.org/package=epiR\");\n\n// Restore the trailing space that used to sit before the removed hyperlink.\nconst p99End = p99.getRange(Word.RangeLocation.end);\np99End.insertText(\" \", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Update bibliography entries: expand abbreviated journal names to their\n# full titles and drop the trailing \"doi:\"/URL hyperlinks (per the\n# \"added updates for phrasing prior to submission\" commit).\n\n$d = $word.ActiveDocument\n\n# --- Reference 1 (paragraph index 90 / Paragraphs.Item(91)): ---\n# Mangtani et al., Clin Infect Dis 2014\n$p = $d.Paragraphs.Item(91)\n$r = $p.Range\n$r.Find.Execute(\"Clin Infect Dis\") | Out-Null\n$wasItalic = $r.Font.Italic\n$r.Text = \"Clinical infectious diseases : an official publication of the Infectious Diseases Society of America\"\n$r.Font.Italic = $wasItalic\n\n$r = $p.Range\n$r.Find.Execute(\" doi:\") | Out-Null\n$r.Delete()\n\n$r = $p.Range\n$r.Find.Execute(\"10.1093/cid/cit790\") | Out-Null\n$r.Delete()\n\n# --- Reference 2 (paragraph index 91 / Paragraphs.Item(92)): ---\n# Abubakar et al., Health Technol Assess 2013\n$p = $d.Paragraphs.Item(92)\n$r = $p.Range\n$r.Find.Execute(\"Health Technol Assess\") | Out-Null\n$wasItalic = $r.Font.Italic\n$r.Text = \"Health technology assessment\"\n$r.Font.Italic = $wasItalic\n\n$r = $p.Range\n$r.Find.Execute(\" doi:\") | Out-Null\n$r.Delete()\n\n$r = $p.Range\n$r.Find.Execute(\"10.3310/hta17370\") | Out-Null\n$r.Delete()\n\n# --- Reference 3 (paragraph index 92 / Paragraphs.Item(93)): ---\n# Zwerling et al., PLoS Med 2011\n$p = $d.Paragraphs.Item(93)\n$r = $p.Range\n$r.Find.Execute(\"PLoS Med\") | Out-Null\n$wasItalic = $r.Font.Italic\n$r.Text = \"PLoS medicine\"\n$r.Font.Italic = $wasItalic\n\n$r = $p.Range\n$r.Find.Execute(\". doi:\") | Out-Null\n$wasBold = $r.Font.Bold\n$r.Text = \":e1001012.\"\nif ($wasBold) { $r.Font.Bold = $wasBold } else { $r.Font.Bold = 0 }\n\n$r = $p.Range\n$r.Find.Execute(\"10.1371/journal.pmed.1001012\") | Out-Null\n$r.Delete()\n\n# --- Reference 4 (paragraph index 93 / Paragraphs.Item(94)): ---\n# Rodrigues et al., Int J Epidemiol 1993\n$p = $d.Paragraphs.Item(94)\n$r = $p.Range\n$r.Find.Execute(\"Int J Epidemiol\") | Out-Null\n$wasItalic = $r.Font.Italic\n$r.Text = \"International journal of epidemiology\"\n$r.Font.Italic = $wasItalic\n\n$r = $p.Range\n$r.Find.Execute(\"http://www.ncbi.nlm.nih.gov/pubmed/8144299\") | Out-Null\n$r.Delete()\n\n# --- Reference 5 (paragraph index 94 / Paragraphs.Item(95)): ---\n# Colditz et al., JAMA 1994 (journal name unchanged)\n$p = $d.Paragraphs.Item(95)\n$r = $p.Range\n$r.Find.Execute(\" doi:\") | Out-Null\n$r.Delete()\n\n$r = $p.Range\n$r.Find.Execute(\"10.1001/jama.1994.03510330076038\") | Out-Null\n$r.Delete()\n\n# --- Reference 6 (paragraph index 95 / Paragraphs.Item(96)): ---\n# Mangtani et al., Clin Infect Dis 2014 (duplicate citation)\n$p = $d.Paragraphs.Item(96)\n$r = $p.Range\n$r.Find.Execute(\"Clin Infect Dis\") | Out-Null\n$wasItalic = $r.Font.Italic\n$r.Text = \"Clinical Infectious Diseases\"\n$r.Font.Italic = $wasItalic\n\n$r = $p.Range\n$r.Find.Execute(\" doi:\") | Out-Null\n$r.Delete()\n\n$r = $p.Range\n$r.Find.Execute(\"10.1093/cid/cit790\") | Out-Null\n$r.Delete()\n\n# --- Reference 7 (paragraph index 96 / Paragraphs.Item(97)): ---\n# Hart & Sutherland, Bull World Health Organ 1972\n$p = $d.Paragraphs.Item(97)\n$r = $p.Range\n$r.Find.Execute(\"Bull World Health Organ\") | Out-Null\n$wasItalic = $r.Font.Italic\n$r.Text = \"The American Statistician\"\n$r.Font.Italic = $wasItalic\n\n$r = $p.Range\n$r.Find.Execute(\" doi:\") | Out-Null\n$r.Delete()\n\n$r = $p.Range\n$r.Find.Execute(\"10.1136/bmj.2.6082.293\") | Out-Null\n$r.Delete()\n\n# --- Reference 10 (paragraph index 99 / Paragraphs.Item(100)): ---\n# Stevenson et al., epiR package 2017\n$p = $d.Paragraphs.Item(100)\n$r = $p.Range\n$r.Find.Execute(\"epiR: Tools for the Analysis of Epidemiological Data\") | Out-Null\n$r.Font.Italic = -1\n\n$r = $p.Range\n$r.Find.Execute(\"https://cran.r-project.org/package=epiR\") | Out-Null\n$r.Delete()\n\n# Restore the trailing space that used to sit before the removed hyperlink.\n$endRange = $p.Range\n$endRange.Collapse(0)\n$endRange.InsertBefore(\" \")\n"}
